$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-41 down to 29-42.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 44489
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112009
$ws.Cells.Item(28, 7).Value = "Acelga"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 1400
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 13).Value = 1450
$ws.Cells.Item(28, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 483
$ws.Cells.Item(28, 17).Value = 3
$ws.Cells.Item(28, 18).Value = "Hortaliza"
